{"js": "// Apply the scripted set of text replacements (title date + the 100 table-cell\n// arithmetic expressions) from the commit's diff, in document order.\nconst replacements = [\n  [\"2023-03-11 Saturday\", \"2023-03-12 Sunday\"],\n  [\"98-61=\", \"98-4=\"],\n  [\"77-50=\", \"57+42=\"],\n  [\"78-21=\", \"16-4=\"],\n  [\"39+0=\", \"86-73=\"],\n  [\"64-24=\", \"82-26=\"],\n  [\"11+66=\", \"8+30=\"],\n  [\"83-79=\", \"51+32=\"],\n  [\"0+63=\", \"24+56=\"],\n  [\"55-49=\", \"62-27=\"],\n  [\"64-21=\", \"89-68=\"],\n  [\"43-18=\", \"68-53=\"],\n  [\"93-15=\", \"18-8=\"],\n  [\"82-37=\", \"5+18=\"],\n  [\"31+22=\", \"45-45=\"],\n  [\"98-51=\", \"81-3=\"],\n  [\"99-50=\", \"5+61=\"],\n  [\"26+16=\", \"70-28=\"],\n  [\"58+20=\", \"72+2=\"],\n  [\"26+39=\", \"4+12=\"],\n  [\"19+68=\", \"57-4=\"],\n  [\"80-59=\", \"61-8=\"],\n  [\"99-11=\", \"66+14=\"],\n  [\"13-2=\", \"89-32=\"],\n  [\"61+37=\", \"3-2=\"],\n  [\"97-41=\", \"74-9=\"],\n  [\"65-63=\", \"48+23=\"],\n  [\"45+17=\", \"44+39=\"],\n  [\"7+62=\", \"19-0=\"],\n  [\"9+29=\", \"65-34=\"],\n  [\"61+32=\", \"81+9=\"],\n  [\"30-30=\", \"56-13=\"],\n  [\"18+45=\", \"16+39=\"],\n  [\"26-10=\", \"66+31=\"],\n  [\"57+7=\", \"78+14=\"],\n  [\"99-73=\", \"95-83=\"],\n  [\"86-9=\", \"54+0=\"],\n  [\"6+52=\", \"71-6=\"],\n  [\"82+1=\", \"4+24=\"],\n  [\"48+9=\", \"47-18=\"],\n  [\"62+2=\", \"30-22=\"],\n  [\"98-29=\", \"78+13=\"],\n  [\"95-47=\", \"34+7=\"],\n  [\"8+19=\", \"30-20=\"],\n  [\"62+23=\", \"14+32=\"],\n  [\"25-2=\", \"12+21=\"],\n  [\"29+37=\", \"60-2=\"],\n  [\"84+14=\", \"58+33=\"],\n  [\"80-75=\", \"61+31=\"],\n  [\"87-12=\", \"78-13=\"],\n  [\"40-0=\", \"73-66=\"],\n  [\"54+16=\", \"78-73=\"],\n  [\"30+40=\", \"71+1=\"],\n  [\"38+57=\", \"41+10=\"],\n  [\"35+13=\", \"59-20=\"],\n  [\"35-7=\", \"91-38=\"],\n  [\"9+63=\", \"30+63=\"],\n  [\"20+3=\", \"80-66=\"],\n  [\"82+11=\", \"56+17=\"],\n  [\"70-50=\", \"44-27=\"],\n  [\"47-23=\", \"66-16=\"],\n  [\"89-65=\", \"73-44=\"],\n  [\"98-93=\", \"2+85=\"],\n  [\"1+57=\", \"59+23=\"],\n  [\"25+38=\", \"94-84=\"],\n  [\"32+55=\", \"50-3=\"],\n  [\"23-6=\", \"98-34=\"],\n  [\"94-80=\", \"90-88=\"],\n  [\"88+6=\", \"40-11=\"],\n  [\"52+30=\", \"13+53=\"],\n  [\"91-16=\", \"85-41=\"],\n  [\"92-34=\", \"37+23=\"],\n  [\"25+23=\", \"36+42=\"],\n  [\"19+29=\", \"85+3=\"],\n  [\"13-7=\", \"24+28=\"],\n  [\"64-1=\", \"20+27=\"],\n  [\"76-26=\", \"15+37=\"],\n  [\"93-32=\", \"53+15=\"],\n  [\"24+13=\", \"58-49=\"],\n  [\"84-70=\", \"13+64=\"],\n  [\"57-19=\", \"89-8=\"],\n  [\"42+48=\", \"36+45=\"],\n  [\"83-12=\", \"11+21=\"],\n  [\"77+5=\", \"33+30=\"],\n  [\"70-2=\", \"84-67=\"],\n  [\"78-36=\", \"78+2=\"],\n  [\"25-9=\", \"8+43=\"],\n  [\"84-59=\", \"53+15=\"],\n  [\"85-3=\", \"52-46=\"],\n  [\"33+29=\", \"5+38=\"],\n  [\"93-18=\", \"27+64=\"],\n  [\"40-35=\", \"37-13=\"],\n  [\"58-3=\", \"27+11=\"],\n  [\"74-22=\", \"83-15=\"],\n  [\"21+68=\", \"31+18=\"],\n  [\"94-38=\", \"5+93=\"],\n  [\"49-16=\", \"42-5=\"],\n  [\"60-46=\", \"19+19=\"],\n  [\"6+26=\", \"87-50=\"],\n  [\"92-29=\", \"38-17=\"],\n  [\"96-68=\", \"33+43=\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load('items');\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(newText, Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n", "ps1": "# Apply the scripted set of text replacements (title date + the 100 table-cell\n# arithmetic expressions) from the commit's diff, in document order.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    ,@(\"2023-03-11 Saturday\", \"2023-03-12 Sunday\")\n    ,@(\"98-61=\", \"98-4=\")\n    ,@(\"77-50=\", \"57+42=\")\n    ,@(\"78-21=\", \"16-4=\")\n    ,@(\"39+0=\", \"86-73=\")\n    ,@(\"64-24=\", \"82-26=\")\n    ,@(\"11+66=\", \"8+30=\")\n    ,@(\"83-79=\", \"51+32=\")\n    ,@(\"0+63=\", \"24+56=\")\n    ,@(\"55-49=\", \"62-27=\")\n    ,@(\"64-21=\", \"89-68=\")\n    ,@(\"43-18=\", \"68-53=\")\n    ,@(\"93-15=\", \"18-8=\")\n    ,@(\"82-37=\", \"5+18=\")\n    ,@(\"31+22=\", \"45-45=\")\n    ,@(\"98-51=\", \"81-3=\")\n    ,@(\"99-50=\", \"5+61=\")\n    ,@(\"26+16=\", \"70-28=\")\n    ,@(\"58+20=\", \"72+2=\")\n    ,@(\"26+39=\", \"4+12=\")\n    ,@(\"19+68=\", \"57-4=\")\n    ,@(\"80-59=\", \"61-8=\")\n    ,@(\"99-11=\", \"66+14=\")\n    ,@(\"13-2=\", \"89-32=\")\n    ,@(\"61+37=\", \"3-2=\")\n    ,@(\"97-41=\", \"74-9=\")\n    ,@(\"65-63=\", \"48+23=\")\n    ,@(\"45+17=\", \"44+39=\")\n    ,@(\"7+62=\", \"19-0=\")\n    ,@(\"9+29=\", \"65-34=\")\n    ,@(\"61+32=\", \"81+9=\")\n    ,@(\"30-30=\", \"56-13=\")\n    ,@(\"18+45=\", \"16+39=\")\n    ,@(\"26-10=\", \"66+31=\")\n    ,@(\"57+7=\", \"78+14=\")\n    ,@(\"99-73=\", \"95-83=\")\n    ,@(\"86-9=\", \"54+0=\")\n    ,@(\"6+52=\", \"71-6=\")\n    ,@(\"82+1=\", \"4+24=\")\n    ,@(\"48+9=\", \"47-18=\")\n    ,@(\"62+2=\", \"30-22=\")\n    ,@(\"98-29=\", \"78+13=\")\n    ,@(\"95-47=\", \"34+7=\")\n    ,@(\"8+19=\", \"30-20=\")\n    ,@(\"62+23=\", \"14+32=\")\n    ,@(\"25-2=\", \"12+21=\")\n    ,@(\"29+37=\", \"60-2=\")\n    ,@(\"84+14=\", \"58+33=\")\n    ,@(\"80-75=\", \"61+31=\")\n    ,@(\"87-12=\", \"78-13=\")\n    ,@(\"40-0=\", \"73-66=\")\n    ,@(\"54+16=\", \"78-73=\")\n    ,@(\"30+40=\", \"71+1=\")\n    ,@(\"38+57=\", \"41+10=\")\n    ,@(\"35+13=\", \"59-20=\")\n    ,@(\"35-7=\", \"91-38=\")\n    ,@(\"9+63=\", \"30+63=\")\n    ,@(\"20+3=\", \"80-66=\")\n    ,@(\"82+11=\", \"56+17=\")\n    ,@(\"70-50=\", \"44-27=\")\n    ,@(\"47-23=\", \"66-16=\")\n    ,@(\"89-65=\", \"73-44=\")\n    ,@(\"98-93=\", \"2+85=\")\n    ,@(\"1+57=\", \"59+23=\")\n    ,@(\"25+38=\", \"94-84=\")\n    ,@(\"32+55=\", \"50-3=\")\n    ,@(\"23-6=\", \"98-34=\")\n    ,@(\"94-80=\", \"90-88=\")\n    ,@(\"88+6=\", \"40-11=\")\n    ,@(\"52+30=\", \"13+53=\")\n    ,@(\"91-16=\", \"85-41=\")\n    ,@(\"92-34=\", \"37+23=\")\n    ,@(\"25+23=\", \"36+42=\")\n    ,@(\"19+29=\", \"85+3=\")\n    ,@(\"13-7=\", \"24+28=\")\n    ,@(\"64-1=\", \"20+27=\")\n    ,@(\"76-26=\", \"15+37=\")\n    ,@(\"93-32=\", \"53+15=\")\n    ,@(\"24+13=\", \"58-49=\")\n    ,@(\"84-70=\", \"13+64=\")\n    ,@(\"57-19=\", \"89-8=\")\n    ,@(\"42+48=\", \"36+45=\")\n    ,@(\"83-12=\", \"11+21=\")\n    ,@(\"77+5=\", \"33+30=\")\n    ,@(\"70-2=\", \"84-67=\")\n    ,@(\"78-36=\", \"78+2=\")\n    ,@(\"25-9=\", \"8+43=\")\n    ,@(\"84-59=\", \"53+15=\")\n    ,@(\"85-3=\", \"52-46=\")\n    ,@(\"33+29=\", \"5+38=\")\n    ,@(\"93-18=\", \"27+64=\")\n    ,@(\"40-35=\", \"37-13=\")\n    ,@(\"58-3=\", \"27+11=\")\n    ,@(\"74-22=\", \"83-15=\")\n    ,@(\"21+68=\", \"31+18=\")\n    ,@(\"94-38=\", \"5+93=\")\n    ,@(\"49-16=\", \"42-5=\")\n    ,@(\"60-46=\", \"19+19=\")\n    ,@(\"6+26=\", \"87-50=\")\n    ,@(\"92-29=\", \"38-17=\")\n    ,@(\"96-68=\", \"33+43=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $true\n    $find.MatchWildcards = $false\n    $find.Execute([ref]$oldText, [ref]$true, [ref]$true, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$false, [ref]$newText, [ref]2) | Out-Null\n}\n"}
